$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (rows 2..31) replacing the previous rows 2..31, and the old
# row 32 is removed entirely (data shifted to the next reporting period).
$data = @(
    @(2, 45969, 8.191000000000001, 0.047),
    @(3, 45969.01041666666, 6.212, 0.498),
    @(4, 45969.02083333334, 9.334, 0),
    @(5, 45969.03125, 8.616, 0),
    @(6, 45969.04166666666, 11.112, 0),
    @(7, 45969.05208333334, 13.144, 0),
    @(8, 45969.0625, 3.548, 1.463),
    @(9, 45969.07291666666, 1.428, 3.758),
    @(10, 45969.08333333334, 0.094, 2.057),
    @(11, 45969.09375, 0.019, 5.8),
    @(12, 45969.10416666666, 1.107, 1.83),
    @(13, 45969.11458333334, 0.731, 5.566),
    @(14, 45969.125, 0.299, 4.991),
    @(15, 45969.13541666666, 0, 9.164),
    @(16, 45969.14583333334, 0, 6.924),
    @(17, 45969.15625, 0, 8.968999999999999),
    @(18, 45969.16666666666, 0, 19.656),
    @(19, 45969.17708333334, 0, 14.908),
    @(20, 45969.1875, 0, 5.321),
    @(21, 45969.19791666666, 0, 5.875),
    @(22, 45969.20833333334, 0, 33.339),
    @(23, 45969.21875, 0, 42.556),
    @(24, 45969.22916666666, 0, 38.377),
    @(25, 45969.23958333334, 0, 13.999),
    @(26, 45969.25, 0, 60.406),
    @(27, 45969.26041666666, 0, 20.337),
    @(28, 45969.27083333334, 0, 13.395),
    @(29, 45969.28125, 0, 12.136),
    @(30, 45969.29166666666, 0, 7.366),
    @(31, 45969.30208333334, 1.866, 16.943)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# The final (32nd) row of the previous dataset no longer exists; remove it
# so the sheet now spans A1:C31.
$ws.Rows.Item(32).Delete()
